$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.5926966292134831
$ws.Range("C2").Value = 0.5527156549520766
$ws.Range("D2").Value = 0.9719101123595506
$ws.Range("E2").Value = 0.7046843177189409
$ws.Range("F2").Value = 0.8439024390243902
$ws.Range("G2").Value = 0.9443627965567919
$ws.Range("H2").Value = 0.7882772938321481
$ws.Range("I2").Value = 519
$ws.Range("J2").Value = 420
$ws.Range("K2").Value = 114
$ws.Range("L2").Value = 15

# ---- Sheet: Classification Report ----
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("B2").Value = 0.8837209302325582
$ws.Range("C2").Value = 0.2134831460674157
$ws.Range("D2").Value = 0.3438914027149321

$ws.Range("B3").Value = 0.5527156549520766
$ws.Range("C3").Value = 0.9719101123595506
$ws.Range("D3").Value = 0.7046843177189409

$ws.Range("B4").Value = 0.5926966292134831
$ws.Range("C4").Value = 0.5926966292134831
$ws.Range("D4").Value = 0.5926966292134831
$ws.Range("E4").Value = 0.5926966292134831

$ws.Range("B5").Value = 0.7182182925923174
$ws.Range("C5").Value = 0.5926966292134832
$ws.Range("D5").Value = 0.5242878602169365

$ws.Range("B6").Value = 0.7182182925923174
$ws.Range("C6").Value = 0.5926966292134831
$ws.Range("D6").Value = 0.5242878602169364

# ---- Sheet: Confusion Matrix ----
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 114
$ws.Range("C2").Value = 420

$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 519
